# Upload new version with timestamp
# - "current balance" values for rows 7 and 8 (H7, H8) are refreshed to "0:0"
# - the footer timestamp (A11) advances from 9:46 AM to 9:47 AM

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "0:0"
$ws.Range("H8").Value = "0:0"
$ws.Range("A11").Value = "Saturday, 24 May, 2025 9:47 AM"
